$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Module Title (row 1, col 2)
$tbl.Cell(1, 2).Range.Text = "Data Preparation & Visualisation, Machine Learning for Data Analytics, Statistics for Data Analytics, Programming for Data Analytics"

# Assessment Title (row 2, col 2)
$tbl.Cell(2, 2).Range.Text = "CA1 50% Integrated Assesment"

# Lecturer Name (row 3, col 2)
$tbl.Cell(3, 2).Range.Text = "David McQuaid, Dr. Muhammad Iqbal, Marina Iantorno, David Gonzalez"

# Assessment Due Date (row 6, col 2) and Date of Submission (row 7, col 2): 11.11.2023 -> 12.11.2023
$d.Content.Find.Execute("11.11.2023", $false, $false, $false, $false, $false,
                         $true, 1, $false, "12.11.2023", 2)
